# Applies the commit: rows 2, 4 and 5 of the "Artfynd" sheet are cyclically
# rotated (new row2 <- old row4, new row4 <- old row5, new row5 <- old row2),
# keeping columns that are identical across the three rows untouched.
#
# Captured "before" values for the columns that actually differ between the
# three rows (everything else - C, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AW, AX, AY - is identical across rows 2/4/5 and is left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- old row 2 (Barrviolspindling / Cortinarius harcynicus) ---
$row2 = @{
    A = 111525963
    B = 85062
    D = "NT"
    E = 249278
    F = "Barrviolspindling"
    G = "Cortinarius harcynicus"
    H = "(Pers.) M.M.Moser"
    Q = 538537.1937094387
    R = 7024283.354085779
}

# --- old row 4 (Knärot / Goodyera repens) ---
$row4 = @{
    A = 111526007
    B = 96348
    D = "VU"
    E = 220787
    F = "Knärot"
    G = "Goodyera repens"
    H = "(L.) R. Br."
    K = "blomning"
    Q = 538522.0815204142
    R = 7024306.075093818
    AC = "Blommande"
}

# --- old row 5 (Rödbrun klubbdyna / Trichoderma nybergianum) ---
$row5 = @{
    A = 111525958
    B = 82949
    D = "NT"
    E = 5589
    F = "Rödbrun klubbdyna"
    G = "Trichoderma nybergianum"
    H = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
    Q = 538520.7165163768
    R = 7024307.405905476
}

# New row 2 gets old row 4's data.
$ws.Range("A2").Value = $row4.A
$ws.Range("B2").Value = $row4.B
$ws.Range("D2").Value = $row4.D
$ws.Range("E2").Value = $row4.E
$ws.Range("F2").Value = $row4.F
$ws.Range("G2").Value = $row4.G
$ws.Range("H2").Value = $row4.H
$ws.Range("K2").Value = $row4.K
$ws.Range("Q2").Value = $row4.Q
$ws.Range("R2").Value = $row4.R
$ws.Range("AC2").Value = $row4.AC

# New row 4 gets old row 5's data; the activity/flowering-related columns
# (J, K, L, N, AC, AF) that only old row 4 had are cleared out.
$ws.Range("A4").Value = $row5.A
$ws.Range("B4").Value = $row5.B
$ws.Range("D4").Value = $row5.D
$ws.Range("E4").Value = $row5.E
$ws.Range("F4").Value = $row5.F
$ws.Range("G4").Value = $row5.G
$ws.Range("H4").Value = $row5.H
$ws.Range("Q4").Value = $row5.Q
$ws.Range("R4").Value = $row5.R
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("AC4").ClearContents()
$ws.Range("AF4").ClearContents()

# New row 5 gets old row 2's data.
$ws.Range("A5").Value = $row2.A
$ws.Range("B5").Value = $row2.B
$ws.Range("D5").Value = $row2.D
$ws.Range("E5").Value = $row2.E
$ws.Range("F5").Value = $row2.F
$ws.Range("G5").Value = $row2.G
$ws.Range("H5").Value = $row2.H
$ws.Range("Q5").Value = $row2.Q
$ws.Range("R5").Value = $row2.R

# New row 2 also gains the (empty) J2/L2/N2/AF2 placeholder cells that came
# along with old row 4's K2="blomning"/AC2="Blommande" values.
$ws.Range("J2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("AF2").Value = ""
